# Apply the "first draft" update to the bpd_rcv results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the "NB" algorithm row (was row 8: A8=6, B8="NB").
#    This shifts the old row 9 (SVM) up to row 8.
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).Delete()

# Fix the index in column A for the row that shifted up (old A=7 -> new A=6)
$ws.Range("A8").Value = 6

# ---------------------------------------------------------------------------
# 2. Rename the "CART" algorithm (row 5) to "DTREE"
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "DTREE"

# ---------------------------------------------------------------------------
# 3. Update header labels (row 1) - existing "Base" columns now denote the
#    mean, and a new "std" column follows each of them.
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "One Year Base mean"
$ws.Range("D1").Value = "One Year Base std"
$ws.Range("E1").Value = "Two Year Base mean"
$ws.Range("F1").Value = "Two Year Base std"
$ws.Range("G1").Value = "Three Year Base mean"
$ws.Range("H1").Value = "Three Year Base std"
$ws.Range("I1").Value = "Five Year Base mean"
$ws.Range("J1").Value = "Five Year Base std"
$ws.Range("K1").Value = "Ten Year Base mean"
$ws.Range("L1").Value = "Ten Year Base std"

# Give the new header cells (H1:L1) the same look (bold/border/center) as
# the rest of the header row.
$ws.Range("B1").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Write out the new mean/std data grid for rows 2-8.
# ---------------------------------------------------------------------------
$data = @(
  @(0.871223292836196,   0.01547207426707886,
    0.8583930869645154,  0.0297326616944997,
    0.8481020877035261,  0.02875459972122795,
    0.8468507751937985,  0.04553773728141705,
    0.8541022469593897,  0.04126511602939186),
  @(0.8712149141181399,  0.01709915973414721,
    0.8672044493473064,  0.02824324318339969,
    0.8594246329038059,  0.03320947084537992,
    0.8639413759689922,  0.04629902670980388,
    0.8692846835703978,  0.03665365690173061),
  @(0.9378759949727693,  0.0139084000407495,
    0.9315407243978674,  0.02030189116296942,
    0.9364049545499951,  0.02186279685885673,
    0.9370578972868217,  0.01873571964885496,
    0.9219954648526075,  0.03418489133270652),
  @(0.786459991621282,   0.04693481522715361,
    0.7716721823864681,  0.03468531481547179,
    0.7845370092897813,  0.03631690083254221,
    0.7698643410852714,  0.04489286189579736,
    0.7527417027417027,  0.03334227679197937),
  @(0.892563887725178,   0.01504154048635178,
    0.8705414598271741,  0.03475021575693216,
    0.8714314254320248,  0.0249456233028538,
    0.8662730135658915,  0.04111317402239689,
    0.8592042877757164,  0.04414232537618996),
  @(0.8809342270632593,  0.01836629953015232,
    0.874632285346571,   0.02932952453737822,
    0.8643442213565079,  0.02230731644275184,
    0.8841751453488372,  0.02959887266797834,
    0.8642238713667284,  0.03630325227844208),
  @(0.8919061583577713,  0.01690275075379055,
    0.8956701599558743,  0.03201505154080678,
    0.8982519228848267,  0.02207250507841166,
    0.9043907461240309,  0.02770665612561962,
    0.8895794681508968,  0.03424640787959609)
)

$cols = @("C","D","E","F","G","H","I","J","K","L")

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $cell = $cols[$j] + $row
        $ws.Range($cell).Value = $values[$j]
    }
}

